# Applies the "Updated documentation and climates" edit:
#  - In the "cases" sheet, the rotation column (F) for every data row changes
#    from "WHEAT.Ble_Dur_1" to "Chickpea.Ghab2"
#  - In the "cases" sheet, the management column (G) for every data row changes
#    from "ROTATION_BLE_IRRIGUE" to "achille_rainfed_3N"
#  - Selection/active cell is updated on three sheets to reflect where the
#    author last clicked while editing.

$wb = $excel.ActiveWorkbook

# -- "cases" sheet: update the rotation/management values for the data rows --
$wsCases = $wb.Worksheets.Item("cases")
$wsCases.Range("F2:F6").Value = """Chickpea.Ghab2"""
$wsCases.Range("G2:G6").Value = """achille_rainfed_3N"""

# -- update active selection on "testble mais poischiche" --
$wsTest = $wb.Worksheets.Item("testble mais poischiche")
$wsTest.Range("E20").Select()

# -- update active selection on "caseswithoutmaize" --
$wsNoMaize = $wb.Worksheets.Item("caseswithoutmaize")
$wsNoMaize.Range("F3").Select()

# -- update active selection on "cases" (also leaves this sheet active/tabSelected) --
$wsCases.Range("G4").Select()
